$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the now-unused "N/A" shared string by replacing G/H (rows with data) "N/A" cells
# with actual numeric values (test results).

$ws.Range("G4").Value = 0.81043026203594604
$ws.Range("H4").Value = 0.76421047685188404

$ws.Range("G5").Value = 0.82003114073988903
$ws.Range("H5").Value = 0.77763542662011798

$ws.Range("G6").Value = 0.81355305833292502
$ws.Range("H6").Value = 0.79161999328722099

$ws.Range("G7").Value = 0.81045561990773995
$ws.Range("H7").Value = 1.21965424699037

$ws.Range("G8").Value = 0.82076074074144201
$ws.Range("H8").Value = 1.1777451317129

$ws.Range("G10").Value = 0.95126683101860998
$ws.Range("H10").Value = 0.95695051331029002

$ws.Range("G11").Value = 0.85050280786957999
$ws.Range("H11").Value = 0.85439606678242996

$ws.Range("G12").Value = 0.84856873194446003
$ws.Range("H12").Value = 0.83986160964504997

$ws.Range("G13").Value = 0.82575110786981998
$ws.Range("H13").Value = 0.84928780297082995

$ws.Range("G14").Value = 0.75164423888894005
$ws.Range("H14").Value = 0.76165194814818005

$ws.Range("G15").Value = 0.73949590138889998
$ws.Range("H15").Value = 0.76850170108025995

$ws.Range("G16").Value = 0.77620039629631998
$ws.Range("H16").Value = 0.79764093827159999

$ws.Range("G17").Value = 0.83138462037036998
$ws.Range("H17").Value = 0.88626506284721995

$ws.Range("G18").Value = 0.75948153888889003
$ws.Range("H18").Value = 0.80326767939815003

$ws.Range("G19").Value = 0.77042453564815006
$ws.Range("H19").Value = 0.81885725563272005

$ws.Range("G21").Value = 0.76798957407326796
$ws.Range("H21").Value = 0.76836267453731

$ws.Range("G22").Value = 0.73466754259277001
$ws.Range("H22").Value = 0.73380146296325199

$ws.Range("G23").Value = 0.72888866666714103
$ws.Range("H23").Value = 0.73438412885805504

$ws.Range("G24").Value = 0.73145039907346099
$ws.Range("H24").Value = 0.76064012407411496

$ws.Range("G25").Value = 0.74139222777833902
$ws.Range("H25").Value = 0.77248586049391099

$ws.Range("G27").Value = 0.78071386851715197
$ws.Range("H27").Value = 0.75688793773109897

$ws.Range("G29").Value = 0.781758357407347
$ws.Range("H29").Value = 0.78634217476849599

$ws.Range("G31").Value = 0.79467081111134596
$ws.Range("H31").Value = 0.79921258194425904

$ws.Range("G32").Value = 0.80448897314784396
$ws.Range("H32").Value = 0.79895860833393095

$ws.Range("G33").Value = 0.78919971018573598
$ws.Range("H33").Value = 0.80330205740715899

$ws.Range("G34").Value = 0.89084253055625595
$ws.Range("H34").Value = 0.93885007314767899

$ws.Range("G35").Value = 0.81282625092663396
$ws.Range("H35").Value = 0.79424108657460402

$ws.Range("G36").Value = 0.79415257962709396
$ws.Range("H36").Value = 0.91232915046283602

# D16 changes from 0 to 7.0E-3, this updates F16 as well via formula recalculation
$ws.Range("D16").Value = 0.0070000000000000001

# B29 was empty, now has a value
$ws.Range("B29").Value = 159

# F28 clears its formula/value entirely (becomes truly blank)
$ws.Range("F28").ClearContents()

# Update active cell selection to H36 as last edited cell
$ws.Range("H36").Select()
